$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.319.65"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "2.048.13"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").Value = "'230.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").Value = "'0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'56.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.03%  "

$ws.Range("D9").Value = "'0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.57%  "

$ws.Range("D10").Value = "'0.0768"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "

$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").Value = "2.350.57"
$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").Value = "'14.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.36%  "

$ws.Range("D14").Value = "'20.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.72%  "

$ws.Range("D15").Value = "'0.752"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.91%  "

$ws.Range("E16").Value = "  -1.99%  "

$ws.Range("D17").Value = "2.049.49"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").Value = "37.318.13"
$ws.Range("E18").Value = "  -1.02%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'5.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.81%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'69.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.57%  "

$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  -1.99%  "

$ws.Range("D22").Value = "'225.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("E24").Value = "  +3.39%  "

$ws.Range("E25").Value = "  -3.81%  "

$ws.Range("D26").Value = "'9.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("D27").Value = "'168.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("E28").Value = "  -5.60%  "

$ws.Range("D29").Value = "'19.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.23%  "

$ws.Range("E30").Value = "  -4.61%  "

$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").Value = "'4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.75%  "

$ws.Range("D33").Value = "'0.0622"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.05%  "

$ws.Range("E34").Value = "  -1.89%  "

$ws.Range("D35").Value = "'2.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.26%  "

$ws.Range("E36").Value = "  -0.37%  "

$ws.Range("D37").Value = "'3.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.24%  "

$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("D39").Value = "'5.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.46%  "

$ws.Range("D40").Value = "'0.0226"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.85%  "

$ws.Range("D41").Value = "'97.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.72%  "

$ws.Range("D42").Value = "1.486.31"
$ws.Range("E42").Value = "  +3.26%  "

$ws.Range("D43").Value = "'2.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").Value = "'0.0949"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("E45").Value = "  +2.90%  "

$ws.Range("D46").Value = "'16.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("D47").Value = "'4.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.78%  "

$ws.Range("E48").Value = "  -3.46%  "

$ws.Range("D49").Value = "'7.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.51%  "

$ws.Range("D50").Value = "'2.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "

$ws.Range("D51").Value = "2.234.33"
$ws.Range("E51").Value = "  -1.53%  "
